$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells (row 1) from _old/_new suffixes to _FV2404/_FV2410
for ($c = 1; $c -le 10; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_old$', '_FV2404')
}
for ($c = 12; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value2 = ($cell.Value2 -replace '_new$', '_FV2410')
}

# Convert the range into an Excel Table (ListObject) named Table1
$range = $ws.Range("A1:U82")
$lo = $ws.ListObjects.Add(1, $range, $null, 1)
$lo.Name = "Table1"

# Freeze the header row (pane split after row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

